$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Nome"
$ws.Range("C1").Value = "VALOR"

# Data values for column C (VALOR)
$valores = @(10, 20, 30, 40, 50, 60, 70, 80, 90, 1, 11, 21, 31, 41, 51)

for ($i = 1; $i -le 15; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = "NOME_$i"
    $ws.Cells.Item($row, 3).Value = $valores[$i - 1]
}

# Column A should use integer number format (not date)
$ws.Range("A2:A16").NumberFormat = "0"

$ws.Range("A1:C16").Select() | Out-Null
